$wb = $excel.ActiveWorkbook

# Map of old sheet name -> new sheet name (pad numeric suffix to 2 digits)
$renames = @{
    "GP1" = "GP01"
    "GP2" = "GP02"
    "BP1" = "BP01"
    "BP2" = "BP02"
    "BP3" = "BP03"
    "BP4" = "BP04"
    "BP5" = "BP05"
    "BP6" = "BP06"
    "BP7" = "BP07"
}

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($renames.ContainsKey($oldName)) {
        $newName = $renames[$oldName]

        # Update the KPI title text in A1, which currently mirrors the old sheet name
        $cell = $ws.Range("A1")
        $oldValue = $cell.Value2
        if ($oldValue -ne $null) {
            $newValue = $oldValue -replace [regex]::Escape($oldName), $newName
            $cell.Value = $newValue
        }

        # Rename the sheet itself
        $ws.Name = $newName
    }
}
